$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.115.24'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.65%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.266.96'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.31%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '305.48'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.24%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '93.26'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.57%  '

$ws.Range("B7").Value = 'BinanceUSD'
$ws.Range("C7").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '69.03'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +6,797.04%  '

$ws.Range("B8").Value = 'XRP'
$ws.Range("C8").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.531'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.17%  '

$ws.Range("B9").Value = 'USDC'
$ws.Range("C9").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.00'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.07%  '

$ws.Range("B10").Value = 'Cardano'
$ws.Range("C10").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.488'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.32%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '32.92'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.36%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0804'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.65%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.112'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.87%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.68'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.48%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.623.87'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.58%  '

$ws.Range("E16").Value = '  +1.73%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.267.94'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.91%  '

$ws.Range("E18").Value = '  +4.08%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '41.980.62'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.51%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.66'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.63%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0919'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.06%  '

$ws.Range("E22").Value = '  +1.29%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '68.22'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.98%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '244.01'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.29%  '

$ws.Range("E25").Value = '  +1.98%  '

$ws.Range("E26").Value = '  +2.82%  '

$ws.Range("E27").Value = '  -0.14%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '23.99'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.57%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.68'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.75%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.09'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.57%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '35.27'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.49%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '159.84'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.58%  '

$ws.Range("E33").Value = '  +3.64%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.00'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.02%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0743'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.02%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.04'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.33%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '17.10'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.16%  '

$ws.Range("E38").Value = '  -1.09%  '

$ws.Range("E39").Value = '  +1.68%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.116'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.86%  '

$ws.Range("E41").Value = '  +0.03%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.03'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.05%  '

$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '19.71'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.26%  '

$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.013.93'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.03%  '

$ws.Range("E45").Value = '  +9.92%  '

$ws.Range("E46").Value = '  +1.76%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.23'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.17%  '

$ws.Range("E48").Value = '  +1.21%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '53.78'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.55%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '72.58'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.82%  '

$ws.Range("E51").Value = '  +0.10%  '
